$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 160.66667
$ws.Range("I4").Value = 91
$ws.Range("K4").Value = 91
$ws.Range("M4").Value = 23
$ws.Range("H51").Value = 1912.2222
$ws.Range("I51").Value = 1940
$ws.Range("J51").Value = 1904.2858
$ws.Range("K51").Value = 1940
$ws.Range("L51").Value = 1904.2858
$ws.Range("M51").Value = -1456
$ws.Range("N51").Value = -2872.2858
$ws.Range("H62").Value = 2851.8333
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 3302.4443
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 3302.4443
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -4550.4443
$ws.Range("H65").Value = 2851.8333
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 3302.4443
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 16512.2215
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -22752.2215
$ws.Range("H100").Value = 2444
$ws.Range("I100").Value = 1600.3636
$ws.Range("K100").Value = 1600.3636
$ws.Range("M100").Value = -1059.3636
$ws.Range("H137").Value = 2766.1345
$ws.Range("I137").Value = 2818.9429
$ws.Range("J137").Value = 2657.4119
$ws.Range("K137").Value = 8456.8287
$ws.Range("L137").Value = 7972.2357
$ws.Range("M137").Value = -5906.8287
$ws.Range("N137").Value = -13072.2357
$ws.Range("H138").Value = 2553.0366
$ws.Range("I138").Value = 1258.6
$ws.Range("J138").Value = 4127.3516
$ws.Range("K138").Value = 3775.8
$ws.Range("L138").Value = 12382.0548
$ws.Range("M138").Value = 1364.2
$ws.Range("N138").Value = -22662.0548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 950.1539
$ws.Range("I2").Value = 915.2
$ws.Range("J2").Value = 1066.6666
$ws.Range("K2").Value = 915.2
$ws.Range("L2").Value = 1066.6666
$ws.Range("M2").Value = -802.2
$ws.Range("N2").Value = -1292.6666
$ws.Range("H61").Value = 1543.6207
$ws.Range("I61").Value = 1413.5385
$ws.Range("K61").Value = 1413.5385
$ws.Range("M61").Value = -1201.5385
$ws.Range("H74").Value = 1313.9117
$ws.Range("I74").Value = 917.1579
$ws.Range("J74").Value = 1816.4667
$ws.Range("K74").Value = 917.1579
$ws.Range("L74").Value = 1816.4667
$ws.Range("M74").Value = -43.15790000000004
$ws.Range("N74").Value = -3564.4667
$ws.Range("H77").Value = 1313.9117
$ws.Range("I77").Value = 917.1579
$ws.Range("J77").Value = 1816.4667
$ws.Range("K77").Value = 4585.7895
$ws.Range("L77").Value = 9082.333500000001
$ws.Range("M77").Value = -217.7894999999999
$ws.Range("N77").Value = -17818.3335
$ws.Range("H116").Value = 950.1539
$ws.Range("I116").Value = 915.2
$ws.Range("J116").Value = 1066.6666
$ws.Range("K116").Value = 915.2
$ws.Range("L116").Value = 1066.6666
$ws.Range("M116").Value = 1378.8
$ws.Range("N116").Value = -5654.6666
$ws.Range("H132").Value = 1897.94
$ws.Range("I132").Value = 1885.6976
$ws.Range("J132").Value = 1973.1428
$ws.Range("K132").Value = 5657.0928
$ws.Range("L132").Value = 5919.428400000001
$ws.Range("M132").Value = -3127.0928
$ws.Range("N132").Value = -10979.4284
$ws.Range("H136").Value = 1543.6207
$ws.Range("I136").Value = 1413.5385
$ws.Range("K136").Value = 4240.6155
$ws.Range("M136").Value = -1690.6155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 950.1539
$ws.Range("I3").Value = 915.2
$ws.Range("J3").Value = 1066.6666
$ws.Range("K3").Value = 915.2
$ws.Range("L3").Value = 1066.6666
$ws.Range("M3").Value = -801.2
$ws.Range("N3").Value = -1294.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 25005
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 50000
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 50000
$ws.Range("M2").Value = 103
$ws.Range("N2").Value = -50226
$ws.Range("H16").Value = 1011.3333
$ws.Range("I16").Value = 1384.2858
$ws.Range("J16").Value = 685
$ws.Range("K16").Value = 1384.2858
$ws.Range("L16").Value = 685
$ws.Range("M16").Value = -1097.2858
$ws.Range("N16").Value = -1259
$ws.Range("H28").Value = 28666.666
$ws.Range("J28").Value = 28666.666
$ws.Range("L28").Value = 28666.666
$ws.Range("N28").Value = -29156.666
$ws.Range("H31").Value = 14506.05
$ws.Range("I31").Value = 18285.316
$ws.Range("K31").Value = 18285.316
$ws.Range("M31").Value = -17990.316
$ws.Range("H34").Value = 14506.05
$ws.Range("I34").Value = 18285.316
$ws.Range("K34").Value = 18285.316
$ws.Range("M34").Value = -18083.316
$ws.Range("H58").Value = 1608.56
$ws.Range("I58").Value = 1676
$ws.Range("J58").Value = 1114
$ws.Range("K58").Value = 1676
$ws.Range("L58").Value = 1114
$ws.Range("M58").Value = -1473
$ws.Range("N58").Value = -1520
$ws.Range("H113").Value = 1011.3333
$ws.Range("I113").Value = 1384.2858
$ws.Range("J113").Value = 685
$ws.Range("K113").Value = 1384.2858
$ws.Range("L113").Value = 685
$ws.Range("M113").Value = 785.7141999999999
$ws.Range("N113").Value = -5025
$ws.Range("H134").Value = 1852.7097
$ws.Range("I134").Value = 1229.2858
$ws.Range("K134").Value = 3687.8574
$ws.Range("M134").Value = -1152.8574
$ws.Range("H136").Value = 1608.56
$ws.Range("I136").Value = 1676
$ws.Range("J136").Value = 1114
$ws.Range("K136").Value = 5028
$ws.Range("L136").Value = 3342
$ws.Range("M136").Value = -2478
$ws.Range("N136").Value = -8442

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 6405.25
$ws.Range("I69").Value = 1853
$ws.Range("K69").Value = 5559
$ws.Range("M69").Value = -4748
$ws.Range("H72").Value = 6405.25
$ws.Range("I72").Value = 1853
$ws.Range("K72").Value = 16677
$ws.Range("M72").Value = -12621
$ws.Range("H113").Value = 998.4838999999999
$ws.Range("I113").Value = 1409.5625
$ws.Range("J113").Value = 560
$ws.Range("K113").Value = 4228.6875
$ws.Range("L113").Value = 1680
$ws.Range("M113").Value = -2058.6875
$ws.Range("N113").Value = -6020
$ws.Range("H130").Value = 2057.5
$ws.Range("J130").Value = 3300
$ws.Range("L130").Value = 9900
$ws.Range("N130").Value = -19940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 759.5217
$ws.Range("I97").Value = 564.6667
$ws.Range("J97").Value = 2805.5
$ws.Range("K97").Value = 564.6667
$ws.Range("L97").Value = 2805.5
$ws.Range("M97").Value = -68.66669999999999
$ws.Range("N97").Value = -3797.5
$ws.Range("H132").Value = 2048.6
$ws.Range("I132").Value = 1908.75
$ws.Range("J132").Value = 4006.5
$ws.Range("K132").Value = 5726.25
$ws.Range("L132").Value = 12019.5
$ws.Range("M132").Value = -3196.25
$ws.Range("N132").Value = -17079.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 873
$ws.Range("I16").Value = 478.6154
$ws.Range("K16").Value = 478.6154
$ws.Range("M16").Value = -308.6154
$ws.Range("H100").Value = 2448.7454
$ws.Range("I100").Value = 1709.9
$ws.Range("J100").Value = 2870.9429
$ws.Range("K100").Value = 1709.9
$ws.Range("L100").Value = 2870.9429
$ws.Range("M100").Value = -1168.9
$ws.Range("N100").Value = -3952.9429
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
$ws.Range("H109").Value = 23900
$ws.Range("J109").Value = 23900
$ws.Range("L109").Value = 23900
$ws.Range("M109").Value = -26674
$ws.Range("H132").Value = 1585.0685
$ws.Range("I132").Value = 1262.4918
$ws.Range("J132").Value = 3224.8333
$ws.Range("K132").Value = 3787.4754
$ws.Range("L132").Value = 9674.499899999999
$ws.Range("M132").Value = -1257.4754
$ws.Range("N132").Value = -14734.4999
$ws.Range("H136").Value = 3570.152
$ws.Range("I136").Value = 2656.68
$ws.Range("J136").Value = 4657.619
$ws.Range("K136").Value = 7970.039999999999
$ws.Range("L136").Value = 13972.857
$ws.Range("M136").Value = -5420.039999999999
$ws.Range("N136").Value = -19072.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 12047.5
$ws.Range("J75").Value = 12047.5
$ws.Range("L75").Value = 12047.5
$ws.Range("N75").Value = -13919.5
$ws.Range("H78").Value = 12047.5
$ws.Range("J78").Value = 12047.5
$ws.Range("L78").Value = 36142.5
$ws.Range("N78").Value = -45502.5
$ws.Range("H99").Value = 27693.777
$ws.Range("I99").Value = 18432
$ws.Range("K99").Value = 18432
$ws.Range("M99").Value = -15437
$ws.Range("H107").Value = 348
$ws.Range("I107").Value = 348
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1044
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 876
$ws.Range("N107").Value = ""
$ws.Range("H136").Value = 568.0909
$ws.Range("I136").Value = 538
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 1614
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = 936
$ws.Range("N136").Value = -8700
